$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.045.01"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").Value = "'2.359.68"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'311.66"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").Value = "'107.65"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  -0.44%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "'40.72"
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("D11").Value = "'0.0913"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("E14").Value = "  -3.56%  "
$ws.Range("D15").Value = "'2.718.70"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "'2.359.27"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "'44.965.21"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").Value = "'14.18"
$ws.Range("E19").Value = "  +9.52%  "
$ws.Range("E20").Value = "  -4.78%  "
$ws.Range("D21").Value = "'0.0000106"
$ws.Range("E21").Value = "  -1.25%  "
$ws.Range("D22").Value = "'72.84"
$ws.Range("E22").Value = "  -1.82%  "
$ws.Range("D23").Value = "'3.51"
$ws.Range("E23").Value = "  +1.59%  "
$ws.Range("D24").Value = "'257.81"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "'2.29"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "'11.05"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  -5.86%  "
$ws.Range("E29").Value = "  +1.43%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0964"
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'22.26"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'36.98"
$ws.Range("E32").Value = "  -5.17%  "
$ws.Range("D33").Value = "'167.69"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").Value = "'3.92"
$ws.Range("E38").Value = "  +3.33%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0351"
$ws.Range("E39").Value = "  -3.42%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.89"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("D42").Value = "'99.53"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "'1.873.87"
$ws.Range("E43").Value = "  +11.93%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").Value = "'69.33"
$ws.Range("E44").Value = "  -3.11%  "
$ws.Range("D45").Value = "'0.228"
$ws.Range("E45").Value = "  -4.58%  "
$ws.Range("E46").Value = "  -0.29%  "
$ws.Range("D47").Value = "'12.74"
$ws.Range("D48").Value = "'80.79"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").Value = "'5.61"
$ws.Range("E49").Value = "  +8.32%  "
$ws.Range("D50").Value = "'110.17"
$ws.Range("E50").Value = "  -3.43%  "
$ws.Range("E51").Value = "  +2.08%  "
